{"js": "// Rename the \"Case QC Passed\" column header to \"Treated QC Passed\", and\n// change every standalone \"Pass\" cell value to \"yes\" (new naming\n// convention for the CNV metrics table). `matchWholeWord` keeps the\n// \"Pass\" search from also hitting the \"Passed\" inside the two QC-header\n// cells (\"Case QC Passed\" / \"Control QC Passed\").\n\nconst body = context.document.body;\n\n// 1) \"Case QC Passed\" -> \"Treated QC Passed\" (table header cell).\nconst headerResults = body.search(\"Case QC Passed\", { matchCase: true });\nheaderResults.load(\"items\");\nawait context.sync();\n\nfor (const r of headerResults.items) {\n  r.insertText(\"Treated QC Passed\", \"Replace\");\n}\n\n// 2) Every whole-word \"Pass\" -> \"yes\" (the per-metric pass/fail cells).\nconst passResults = body.search(\"Pass\", {\n  matchCase: true,\n  matchWholeWord: true,\n});\npassResults.load(\"items\");\nawait context.sync();\n\nfor (const r of passResults.items) {\n  r.insertText(\"yes\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Rename the \"Case QC Passed\" column header to \"Treated QC Passed\", and\n# change every standalone \"Pass\" cell value to \"yes\" (new naming\n# convention for the CNV metrics table). MatchWholeWord keeps the\n# \"Pass\" search from also hitting the \"Passed\" inside the two QC-header\n# cells (\"Case QC Passed\" / \"Control QC Passed\").\n\n$d = $word.ActiveDocument\n\n# 1) \"Case QC Passed\" -> \"Treated QC Passed\" (table header cell).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Case QC Passed\"\n$find.Replacement.Text = \"Treated QC Passed\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2) Every whole-word \"Pass\" -> \"yes\" (the per-metric pass/fail cells).\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Pass\"\n$find2.Replacement.Text = \"yes\"\n$find2.Forward = $true\n$find2.Wrap = 1\n$find2.Format = $false\n$find2.MatchCase = $true\n$find2.MatchWholeWord = $true\n$find2.Execute($find2.Text, $find2.MatchCase, $find2.MatchWholeWord, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
